$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 87981.0709163148
$ws.Range("D2").Value = 10137753.70137369
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 10225982.1040372

$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 4.371470058157054

$ws.Range("B4").Value = 0.7287194209349384
$ws.Range("C4").Value = 0.3375848360084654
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 7.700638116232206

$ws.Range("B5").Value = 0.1554434735375247
$ws.Range("C5").Value = 0.3375848360084654
$ws.Range("D5").Value = 16.98373111632243
$ws.Range("E5").Value = 6.48142807727062
$ws.Range("G5").Value = 23.95818750313904
